$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (MuSCs/Resolving-Mac target-cluster duplicate rows removed by the
# updated TPM script output) before touching the remaining values so the row
# numbers used below stay aligned with rows 2-5.
$ws.Rows("6:9").Delete() | Out-Null

# Row 2: Sending cluster ECs -> Target cluster changes from ECs to FAPs, and the
# TPM-derived numeric columns are recalculated with the new values.
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.555934
$ws.Range("H2").Value = 1.667802
$ws.Range("I2").Value = 0.005745252779589096
$ws.Range("J2").Value = 0.005745252779589094
$ws.Range("M2").Value = 0.01688366666666667
$ws.Range("N2").Value = 0.050651
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.009386204344666669
$ws.Range("R2").Value = 0.084475839102
$ws.Range("S2").Value = 0.005745252779589096
$ws.Range("T2").Value = 0.005745252779589094

# Row 3: Sending cluster changes from ECs to FAPs; target cluster stays FAPs.
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 79.68771233333334
$ws.Range("H3").Value = 239.063137
$ws.Range("I3").Value = 0.823525905561055
$ws.Range("J3").Value = 0.823525905561055
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.345420772465223
$ws.Range("R3").Value = 12.108786952187
$ws.Range("S3").Value = 0.823525905561055
$ws.Range("T3").Value = 0.823525905561055

# Row 4: Sending cluster changes from FAPs to MuSCs; target cluster changes from
# ECs to FAPs.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.3446996666666666
$ws.Range("H4").Value = 1.034099
$ws.Range("I4").Value = 0.003562269474506148
$ws.Range("J4").Value = 0.003562269474506148
$ws.Range("M4").Value = 0.01688366666666667
$ws.Range("N4").Value = 0.050651
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.005819794272111111
$ws.Range("R4").Value = 0.052378148449
$ws.Range("S4").Value = 0.003562269474506148
$ws.Range("T4").Value = 0.003562269474506148

# Row 5: Sending cluster changes from FAPs to Resolving-Mac; target cluster stays FAPs.
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 16.17571666666667
$ws.Range("H5").Value = 48.52715
$ws.Range("I5").Value = 0.1671665721848498
$ws.Range("J5").Value = 0.1671665721848498
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.2731054082944445
$ws.Range("R5").Value = 2.45794867465
$ws.Range("S5").Value = 0.1671665721848498
$ws.Range("T5").Value = 0.1671665721848498
